$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking row): Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total row): Total right marks 130 -> 104, Wrong marking -1 -> -2
$ws.Range("B12").Value = 104
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "102 / 112"
